# B1--and-B2-PowerPoint.pptx edit
#
# 1) The table on slide 5 gets a new table style id.
# 2) The deck's theme ("Integral" / Red Violet, carried by the slide master)
#    and the notes-theme ("Office Theme", carried by the notes master) trade
#    their color schemes: the slides end up using the stock Office palette
#    while the notes master ends up with the old Red-Violet palette. (The
#    font scheme and format scheme are already identical between the two
#    themes, so only the twelve theme colors - and the table style - are
#    actually different content.)

function HexToRgbInt($hex) {
    $rr = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $gg = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $bb = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $rr + ($gg * 256) + ($bb * 65536)
}

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 -------------------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{948BE4C6-74FA-406C-8F05-35763302CE71}")

# --- 2. Swap the theme color palettes ------------------------------------
# Colors() index map: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6,
#                      11 hlink, 12 folHlink
$slideTcs = $p.SlideMaster.Theme.ThemeColorScheme

# New ("Office Theme") palette for the slide master's theme:
$slideTcs.Colors(1).RGB  = HexToRgbInt "000000"   # dk1
$slideTcs.Colors(2).RGB  = HexToRgbInt "FFFFFF"   # lt1
$slideTcs.Colors(3).RGB  = HexToRgbInt "44546A"   # dk2
$slideTcs.Colors(4).RGB  = HexToRgbInt "E7E6E6"   # lt2
$slideTcs.Colors(5).RGB  = HexToRgbInt "5B9BD5"   # accent1
$slideTcs.Colors(6).RGB  = HexToRgbInt "ED7D31"   # accent2
$slideTcs.Colors(7).RGB  = HexToRgbInt "A5A5A5"   # accent3
$slideTcs.Colors(8).RGB  = HexToRgbInt "FFC000"   # accent4
$slideTcs.Colors(9).RGB  = HexToRgbInt "4472C4"   # accent5
$slideTcs.Colors(10).RGB = HexToRgbInt "70AD47"   # accent6
$slideTcs.Colors(11).RGB = HexToRgbInt "0563C1"   # hlink
$slideTcs.Colors(12).RGB = HexToRgbInt "954F72"   # folHlink
